$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The "支付协议" (payment agreement) row's access Method changed from POST to GET.
$ws.Range("C4").Value = "GET"

# 2. Add a new row (row 12) describing the "预订协议" (booking/reservation agreement)
#    security-H5 page, following the same layout/pattern as the existing rows.

# First give the new row B12 a hyperlink (reusing the same target used by the
# "账单计划" row right above it, row 11) before formatting is copied over, so the
# paste-formats step below can restyle the cell with the table's normal bordered
# hyperlink look instead of Excel's default "just-added-hyperlink" style.
$ws.Hyperlinks.Add($ws.Range("B12"), "http://api.mogoroom.com/h5/budget", "", "", "http://api.mogoroom.com/h5/budget")

# Clone the formatting (borders/fonts/fills) of row 11 onto row 12.
$ws.Range("A11:I11").Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)

# I3/I4 (and by extension C-column cells sharing that look) use a narrower
# left/right-only border style; the new row's Method (C12) and HTML5-file (I12)
# columns use that same look, so copy it in specifically for those two cells.
$ws.Range("I3").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("I12").PasteSpecial(-4122)

# A12 (page name) keeps the workbook's default (unstyled) look.
$ws.Range("A12").ClearFormats()

# Now fill in the row's actual content.
$ws.Range("A12").Value = "预订协议"
$ws.Range("B12").Value = "http://{domain}/securityH5/book"
$ws.Range("C12").Value = "POST"
$ws.Range("D12").Value = '{"head":{"userId":"","token":"","os":"","osVersion":"","appVersion":"","model":"","uuid":"","channel":"租客app","key":""}}'
$ws.Range("E12").Value = "无附加参数"
$ws.Range("F12").Value = "html文本"
$ws.Range("G12").Value = "否"
$ws.Range("H12").Value = "是"
$ws.Range("I12").Value = "/pages/terms/book.jsp"

# 3. Leave the cursor roughly where the author left it when saving.
$ws.Range("B32").Select()
